$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, pushing the existing rows 2-4 down to 3-5.
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row 2 with the column "slug" identifiers that
# correspond to each header in row 1. These values allow metadata columns to
# be related to each other (SKOS hierarchy support).
$ws.Range("A2").Value = "mes-codigo"
$ws.Range("B2").Value = "numero-de-contratos"
$ws.Range("C2").Value = "provincia-codigo"
$ws.Range("D2").Value = "provincia-nombre"
$ws.Range("E2").Value = "sexo"
$ws.Range("F2").Value = "mes-y-ano"

# Match the styling used by the rest of the sheet's data rows.
$ws.Range("A2:F2").Font.Name = $ws.Range("A3:F3").Font.Name
$ws.Range("A2:F2").Font.Size = $ws.Range("A3:F3").Font.Size
